$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the measured data in rows 3-5 (G..K); columns I, L, M are formulas
# and recalculate automatically.
$ws.Range("G3:G5").Value = 1752
$ws.Range("H3:H5").Value = 1672
$ws.Range("J3:J5").Value = 185.63
$ws.Range("K3:K5").Value = 189

# Move the active selection to L7, matching the saved workbook state.
$ws.Range("L7").Select()
